# Rename ObjTables metadata attributes from UpperCamelCase to lowerCamelCase
$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("!!Data repo metadata")
$wsData.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$wsData.Range("A2").Value = "!!ObjTables type='Data' id='DataRepoMetadata'"

$wsSchema = $wb.Worksheets.Item("!!Schema repo metadata")
$wsSchema.Range("A1").Value = "!!ObjTables type='Data' id='SchemaRepoMetadata'"
